$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "é maior"
$ws.Range("E4").Value = "é maior"
$ws.Range("E3").Value = "é menor"
$ws.Range("E5").Value = "é menor"
